# Apply the fixed standard-curve values to Sheet1 and update the selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the standard curve intercept values (D2, D3).
$ws.Range("D2").Value = 43.146
$ws.Range("D3").Value = 43.487

# Activate the sheet and set the selection to A2:E7 with A2 as the active cell.
$ws.Activate()
$ws.Range("A2:E7").Select()
